$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.130.88'
$ws.Range("E2").Value = '  +3.20%  '
$ws.Range("D3").Value = '3.041.69'
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''595.78'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = '''154.58'
$ws.Range("E6").Value = '  +8.18%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.040.89'
$ws.Range("E8").Value = '  +2.22%  '
$ws.Range("D9").Value = '''0.517'
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").Value = '''6.92'
$ws.Range("E10").Value = '  +14.95%  '
$ws.Range("E11").Value = '  +4.70%  '
$ws.Range("E12").Value = '  +2.74%  '
$ws.Range("E13").Value = '  +3.74%  '
$ws.Range("D14").Value = '''36.07'
$ws.Range("E14").Value = '  +5.40%  '
$ws.Range("E15").Value = '  +2.16%  '
$ws.Range("D16").Value = '3.546.07'
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("D17").Value = '''7.12'
$ws.Range("E17").Value = '  +3.05%  '
$ws.Range("D18").Value = '63.075.73'
$ws.Range("E18").Value = '  +3.03%  '
$ws.Range("D19").Value = '3.045.40'
$ws.Range("E19").Value = '  +1.94%  '
$ws.Range("D20").Value = '''455.74'
$ws.Range("E20").Value = '  +1.76%  '
$ws.Range("E21").Value = '  +3.05%  '
$ws.Range("E22").Value = '  +2.89%  '
$ws.Range("E23").Value = '  +3.70%  '
$ws.Range("E24").Value = '  +2.16%  '
$ws.Range("D25").Value = '''11.32'
$ws.Range("E25").Value = '  +6.10%  '
$ws.Range("D26").Value = '''2.31'
$ws.Range("E26").Value = '  +5.35%  '
$ws.Range("E27").Value = '  +4.30%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +4.62%  '
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("E31").Value = '  +9.73%  '
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").Value = '''27.74'
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("E35").Value = '  +5.24%  '
$ws.Range("D36").Value = '''1.04'
$ws.Range("E36").Value = '  +3.00%  '
$ws.Range("E37").Value = '  +3.72%  '
$ws.Range("D38").Value = '''3.19'
$ws.Range("E38").Value = '  +11.86%  '
$ws.Range("E39").Value = '  +8.03%  '
$ws.Range("D40").Value = '''2.11'
$ws.Range("E40").Value = '  +3.11%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '''9.14'
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("E43").Value = '  +13.61%  '
$ws.Range("D44").Value = '''43.60'
$ws.Range("E44").Value = '  +10.38%  '
$ws.Range("D45").Value = '''393.81'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").Value = '''0.0362'
$ws.Range("E46").Value = '  +3.39%  '
$ws.Range("D47").Value = '2.727.62'
$ws.Range("E47").Value = '  +1.64%  '
$ws.Range("D48").Value = '''132.39'
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("E50").Value = '  +7.40%  '
$ws.Range("D51").Value = '''24.61'
$ws.Range("E51").Value = '  +5.50%  '
